# CU-35_ConsultarReporteDeVentas / Descripcion.docx
#
# "Flujo normal" step that documents which entity fields the system
# reads from the database. The sentence currently reads:
#
#   ... VENTA (NoVenta, FechaRegistro, Total) y CAJA (NoCaja, Cajero) (EX-01), ...
#
# and must become:
#
#   ... VENTA (noVenta, fechaRegistro, total), CAJA (noCaja) y EMPLEADO
#       (nombre, apellidoPaterno) (EX-01), ...
#
# i.e. the VENTA/CAJA field lists are lower-cased, the CAJA clause drops
# "Cajero" (moved out to its own EMPLEADO entity) and a new EMPLEADO
# (nombre, apellidoPaterno) clause is added right after CAJA (noCaja).
# The leading "... de la " lead-in and the trailing " (EX-01), muestra..."
# tail are untouched.

$d = $word.ActiveDocument

$oldPhrase = "VENTA (NoVenta, FechaRegistro, Total) y CAJA (NoCaja, Cajero)"
$newPhrase = "VENTA (noVenta, fechaRegistro, total), CAJA (noCaja) y EMPLEADO (nombre, apellidoPaterno)"

# Locate the phrase without mutating the document (Find on a Duplicate of
# the whole story range), then rewrite just that span in place. Using a
# freshly-minted Range (rather than the Find range itself) for the
# mutation is important: this host's COM shim only commits edits made
# through a live Range/Selection object.
$search = $d.Content.Duplicate
$found = $search.Find.Execute($oldPhrase, $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

if ($found) {
    $target = $d.Range($search.Start, $search.End)
    $target.Text = $newPhrase
} else {
    throw "Could not find expected phrase '$oldPhrase' in the document."
}
